# Danh_Muc_Nhan_Vien_Import.xlsx -- "fix lay nhat ky mau hoa don va file excel"
#
# 1) Rename the worksheet KETOAN-BACHKHOA -> HOADON-BACHKHOA
# 2) Move the selection from K5 to the full column A (A1:A1048576)
#    (the last-used-cell/zoom bookkeeping the diff also touches is
#    Excel's own autosave-view metadata, not user-addressable state)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KETOAN-BACHKHOA")

$ws.Name = "HOADON-BACHKHOA"

[void]$ws.Range("A1:A1048576").Select()
